$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.544.33'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '3.063.77'
$ws.Range("E3").Value = '  +2.68%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '386.07'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.06'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.40%  '
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.584'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.30%  '
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").Value = '3.552.81'
$ws.Range("E13").Value = '  +2.89%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.65'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.06%  '
$ws.Range("E15").Value = '  -0.21%  '
$ws.Range("D16").Value = '3.064.67'
$ws.Range("E16").Value = '  +2.48%  '
$ws.Range("E17").Value = '  -2.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.71'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -3.56%  '
$ws.Range("D19").Value = '51.621.13'
$ws.Range("E19").Value = '  +0.31%  '
$ws.Range("E20").Value = '  +2.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.46'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.20%  '
$ws.Range("E22").Value = '  +0.69%  '
$ws.Range("E23").Value = '  -0.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.99'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.15'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.18'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +4.24%  '
$ws.Range("E27").Value = '  +3.13%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.170'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.56%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.26'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.35%  '
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("E31").Value = '  -2.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.26'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.81'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.57%  '
$ws.Range("E34").Value = '  +0.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.09'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.91%  '
$ws.Range("E36").Value = '  +1.95%  '
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.32'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.54%  '
$ws.Range("E39").Value = '  +7.90%  '
$ws.Range("E40").Value = '  +1.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '16.90'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.68%  '
$ws.Range("E42").Value = '  +0.39%  '
$ws.Range("E43").Value = '  -0.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '125.44'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.75'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.02'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.08'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.44'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +2.09%  '
$ws.Range("D49").Value = '2.037.37'
$ws.Range("E49").Value = '  +0.23%  '
$ws.Range("D50").Value = '3.365.26'
$ws.Range("E50").Value = '  +2.52%  '
$ws.Range("E51").Value = '  +7.35%  '
